$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '66.390.68'
$cell.Style = $origStyle
$cell = $ws.Range("E2")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -0.48%  '
$cell.Style = $origStyle
$cell = $ws.Range("D3")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.461.96'
$cell.Style = $origStyle
$cell = $ws.Range("E3")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -1.58%  '
$cell.Style = $origStyle
$cell = $ws.Range("E4")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +0.03%  '
$cell.Style = $origStyle
$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '583.64'
$cell.Style = $origStyle
$cell = $ws.Range("E5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -0.17%  '
$cell.Style = $origStyle
$cell = $ws.Range("D6")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '177.75'
$cell.Style = $origStyle
$cell = $ws.Range("E6")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +1.11%  '
$cell.Style = $origStyle
$cell = $ws.Range("D7")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.629'
$cell.Style = $origStyle
$cell = $ws.Range("E7")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +5.42%  '
$cell.Style = $origStyle
$cell = $ws.Range("D9")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.461.57'
$cell.Style = $origStyle
$cell = $ws.Range("E9")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -1.45%  '
$cell.Style = $origStyle
$cell = $ws.Range("D10")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.133'
$cell.Style = $origStyle
$cell = $ws.Range("E10")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -1.05%  '
$cell.Style = $origStyle
$cell = $ws.Range("D11")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '6.97'
$cell.Style = $origStyle
$cell = $ws.Range("E11")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +0.77%  '
$cell.Style = $origStyle
$cell = $ws.Range("E12")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -1.52%  '
$cell.Style = $origStyle
$cell = $ws.Range("D13")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '4.062.13'
$cell.Style = $origStyle
$cell = $ws.Range("E13")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -1.51%  '
$cell.Style = $origStyle
$cell = $ws.Range("E14")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +1.35%  '
$cell.Style = $origStyle
$cell = $ws.Range("D15")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '30.21'
$cell.Style = $origStyle
$cell = $ws.Range("E15")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -1.11%  '
$cell.Style = $origStyle
$cell = $ws.Range("D16")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '66.301.26'
$cell.Style = $origStyle
$cell = $ws.Range("E16")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -0.59%  '
$cell.Style = $origStyle
$cell = $ws.Range("E17")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -0.65%  '
$cell.Style = $origStyle
$cell = $ws.Range("D18")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.448.86'
$cell.Style = $origStyle
$cell = $ws.Range("E18")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -1.79%  '
$cell.Style = $origStyle
$cell = $ws.Range("D19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '5.98'
$cell.Style = $origStyle
$cell = $ws.Range("E19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -1.29%  '
$cell.Style = $origStyle
$cell = $ws.Range("D20")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '13.86'
$cell.Style = $origStyle
$cell = $ws.Range("E20")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -1.14%  '
$cell.Style = $origStyle
$cell = $ws.Range("D21")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '371.37'
$cell.Style = $origStyle
$cell = $ws.Range("E21")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -2.84%  '
$cell.Style = $origStyle
$cell = $ws.Range("E22")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -3.16%  '
$cell.Style = $origStyle
$cell = $ws.Range("D23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '73.25'
$cell.Style = $origStyle
$cell = $ws.Range("E23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +1.22%  '
$cell.Style = $origStyle
$cell = $ws.Range("D24")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.997'
$cell.Style = $origStyle
$cell = $ws.Range("E24")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -0.40%  '
$cell.Style = $origStyle
$cell = $ws.Range("D25")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.536'
$cell.Style = $origStyle
$cell = $ws.Range("E25")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -2.76%  '
$cell.Style = $origStyle
$cell = $ws.Range("D26")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0000126'
$cell.Style = $origStyle
$cell = $ws.Range("E26")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +3.98%  '
$cell.Style = $origStyle
$cell = $ws.Range("D27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '10.04'
$cell.Style = $origStyle
$cell = $ws.Range("E27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +1.38%  '
$cell.Style = $origStyle
$cell = $ws.Range("D28")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.178'
$cell.Style = $origStyle
$cell = $ws.Range("E28")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +2.80%  '
$cell.Style = $origStyle
$cell = $ws.Range("E29")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +0.08%  '
$cell.Style = $origStyle
$cell = $ws.Range("D30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '5.97'
$cell.Style = $origStyle
$cell = $ws.Range("E30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +0.61%  '
$cell.Style = $origStyle
$cell = $ws.Range("E31")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -1.22%  '
$cell.Style = $origStyle
$cell = $ws.Range("D32")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '23.68'
$cell.Style = $origStyle
$cell = $ws.Range("E32")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -3.81%  '
$cell.Style = $origStyle
$cell = $ws.Range("E33")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +0.01%  '
$cell.Style = $origStyle
$cell = $ws.Range("E34")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -2.75%  '
$cell.Style = $origStyle
$cell = $ws.Range("E35")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -5.91%  '
$cell.Style = $origStyle
$cell = $ws.Range("E36")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -0.67%  '
$cell.Style = $origStyle
$cell = $ws.Range("D37")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '161.29'
$cell.Style = $origStyle
$cell = $ws.Range("E37")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -0.17%  '
$cell.Style = $origStyle
$cell = $ws.Range("D38")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.887'
$cell.Style = $origStyle
$cell = $ws.Range("E38")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -1.28%  '
$cell.Style = $origStyle
$cell = $ws.Range("D39")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '27.92'
$cell.Style = $origStyle
$cell = $ws.Range("E39")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -7.48%  '
$cell.Style = $origStyle
$cell = $ws.Range("E40")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +0.83%  '
$cell.Style = $origStyle
$cell = $ws.Range("D41")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.811.71'
$cell.Style = $origStyle
$cell = $ws.Range("E41")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +3.02%  '
$cell.Style = $origStyle
$cell = $ws.Range("E42")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +0.40%  '
$cell.Style = $origStyle
$cell = $ws.Range("E43")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +1.39%  '
$cell.Style = $origStyle
$cell = $ws.Range("D44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '6.46'
$cell.Style = $origStyle
$cell = $ws.Range("E44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -0.72%  '
$cell.Style = $origStyle
$cell = $ws.Range("D45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0694'
$cell.Style = $origStyle
$cell = $ws.Range("E45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -1.36%  '
$cell.Style = $origStyle
$cell = $ws.Range("D46")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '25.25'
$cell.Style = $origStyle
$cell = $ws.Range("E46")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +0.63%  '
$cell.Style = $origStyle
$cell = $ws.Range("D47")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '341.72'
$cell.Style = $origStyle
$cell = $ws.Range("E47")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +4.94%  '
$cell.Style = $origStyle
$cell = $ws.Range("D48")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '40.03'
$cell.Style = $origStyle
$cell = $ws.Range("E48")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -1.77%  '
$cell.Style = $origStyle
$cell = $ws.Range("D49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0291'
$cell.Style = $origStyle
$cell = $ws.Range("E49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -0.77%  '
$cell.Style = $origStyle
$cell = $ws.Range("E50")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +2.76%  '
$cell.Style = $origStyle
$cell = $ws.Range("D51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '31.70'
$cell.Style = $origStyle
$cell = $ws.Range("E51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +2.23%  '
$cell.Style = $origStyle
